$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows at the top of the data block (rows 435-437), shifting
# existing rows 435:560 down to 438:563. Excel will auto-update the sheet
# dimension and preserve all existing content/formatting of shifted rows.
$ws.Rows("435:437").Insert()

# --- Fill in the new row 435 (Pinton) ---
$ws.Cells.Item(435,1).Value = 8
$ws.Cells.Item(435,2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(435,3).Value = "Coquimbo"
$ws.Cells.Item(435,4).Value = 44627
$ws.Cells.Item(435,5).Value = 4
$ws.Cells.Item(435,6).Value = "Fruta"
$ws.Cells.Item(435,7).Value = 100108
$ws.Cells.Item(435,8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(435,9).Value = 100108006
$ws.Cells.Item(435,10).Value = "Plátano"
$ws.Cells.Item(435,11).Value = "Sin especificar"
$ws.Cells.Item(435,12).Value = "Pintón"
$ws.Cells.Item(435,13).Value = 80
$ws.Cells.Item(435,14).Value = 19000
$ws.Cells.Item(435,15).Value = 19000
$ws.Cells.Item(435,16).Value = 19000
$ws.Cells.Item(435,17).Value = "`$/caja 20 kilos"
$ws.Cells.Item(435,18).Value = "Ecuador"
$ws.Cells.Item(435,19).Value = 950
$ws.Cells.Item(435,20).Value = 20

# --- Fill in the new row 436 (Primera Maduro) ---
$ws.Cells.Item(436,1).Value = 8
$ws.Cells.Item(436,2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(436,3).Value = "Coquimbo"
$ws.Cells.Item(436,4).Value = 44627
$ws.Cells.Item(436,5).Value = 4
$ws.Cells.Item(436,6).Value = "Fruta"
$ws.Cells.Item(436,7).Value = 100108
$ws.Cells.Item(436,8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(436,9).Value = 100108006
$ws.Cells.Item(436,10).Value = "Plátano"
$ws.Cells.Item(436,11).Value = "Sin especificar"
$ws.Cells.Item(436,12).Value = "Primera Maduro"
$ws.Cells.Item(436,13).Value = 120
$ws.Cells.Item(436,14).Value = 20000
$ws.Cells.Item(436,15).Value = 20000
$ws.Cells.Item(436,16).Value = 20000
$ws.Cells.Item(436,17).Value = "`$/caja 20 kilos"
$ws.Cells.Item(436,18).Value = "Ecuador"
$ws.Cells.Item(436,19).Value = 1000
$ws.Cells.Item(436,20).Value = 20

# --- Fill in the new row 437 (Primera Pintón) ---
$ws.Cells.Item(437,1).Value = 8
$ws.Cells.Item(437,2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(437,3).Value = "Coquimbo"
$ws.Cells.Item(437,4).Value = 44627
$ws.Cells.Item(437,5).Value = 4
$ws.Cells.Item(437,6).Value = "Fruta"
$ws.Cells.Item(437,7).Value = 100108
$ws.Cells.Item(437,8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(437,9).Value = 100108006
$ws.Cells.Item(437,10).Value = "Plátano"
$ws.Cells.Item(437,11).Value = "Sin especificar"
$ws.Cells.Item(437,12).Value = "Primera Pintón"
$ws.Cells.Item(437,13).Value = 128
$ws.Cells.Item(437,14).Value = 21000
$ws.Cells.Item(437,15).Value = 21000
$ws.Cells.Item(437,16).Value = 21000
$ws.Cells.Item(437,17).Value = "`$/caja 20 kilos"
$ws.Cells.Item(437,18).Value = "Ecuador"
$ws.Cells.Item(437,19).Value = 1050
$ws.Cells.Item(437,20).Value = 20

# Apply the date style (numFmt 165) used elsewhere in column D to the new D cells
$ws.Cells.Item(435,4).NumberFormat = $ws.Cells.Item(438,4).NumberFormat
$ws.Cells.Item(436,4).NumberFormat = $ws.Cells.Item(438,4).NumberFormat
$ws.Cells.Item(437,4).NumberFormat = $ws.Cells.Item(438,4).NumberFormat
